# Apply cryptocurrency price/volume updates as described in the commit diff.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Row 47 and 48 swap B (Coin) / C (Link) / D (Price) / E (Volume) content between
# NEARProtocol and WEMIXToken, plus all other numeric price/volume refreshes below.

$ws.Range("D2").Value = '28.003.39'
$ws.Range("E2").Value = '  -1.95%  '
$ws.Range("D3").Value = '1.828.72'
$ws.Range("E3").Value = '  -1.05%  '
$ws.Range("D5").NumberFormat = "@"
$ws.Range("D5").Value = '325.15'
$ws.Range("E5").Value = '  -2.51%  '
$ws.Range("E6").Value = '  -0.06%  '
$ws.Range("D7").NumberFormat = "@"
$ws.Range("D7").Value = '0.4652'
$ws.Range("E7").Value = '  -0.18%  '
$ws.Range("D8").NumberFormat = "@"
$ws.Range("D8").Value = '0.3864'
$ws.Range("E8").Value = '  -1.49%  '
$ws.Range("D9").NumberFormat = "@"
$ws.Range("D9").Value = '0.07855'
$ws.Range("E9").Value = '  -0.71%  '
$ws.Range("D10").NumberFormat = "@"
$ws.Range("D10").Value = '0.9583'
$ws.Range("D11").NumberFormat = "@"
$ws.Range("D11").Value = '21.85'
$ws.Range("E11").Value = '  -1.61%  '
$ws.Range("D12").Value = '1.840.99'
$ws.Range("E12").Value = '  -7.02%  '
$ws.Range("D13").NumberFormat = "@"
$ws.Range("D13").Value = '5.673'
$ws.Range("E13").Value = '  -3.02%  '
$ws.Range("D14").NumberFormat = "@"
$ws.Range("D14").Value = '6.887'
$ws.Range("E14").Value = '  -1.97%  '
$ws.Range("E15").Value = '  -0.53%  '
$ws.Range("D16").NumberFormat = "@"
$ws.Range("D16").Value = '87.14'
$ws.Range("E16").Value = '  -0.68%  '
$ws.Range("D17").NumberFormat = "@"
$ws.Range("D17").Value = '1.002'
$ws.Range("E17").Value = '  -0.11%  '
$ws.Range("D18").NumberFormat = "@"
$ws.Range("D18").Value = '0.000009916'
$ws.Range("E18").Value = '  -1.50%  '
$ws.Range("D19").NumberFormat = "@"
$ws.Range("D19").Value = '16.59'
$ws.Range("E19").Value = '  -3.02%  '
$ws.Range("D20").NumberFormat = "@"
$ws.Range("D20").Value = '0.9999'
$ws.Range("E20").Value = '  -0.19%  '
$ws.Range("D21").Value = '28.025.16'
$ws.Range("E21").Value = '  -1.97%  '
$ws.Range("D22").NumberFormat = "@"
$ws.Range("D22").Value = '5.312'
$ws.Range("E22").Value = '  -1.67%  '
$ws.Range("D23").NumberFormat = "@"
$ws.Range("D23").Value = '10.96'
$ws.Range("E23").Value = '  -2.73%  '
$ws.Range("E24").Value = '  -1.74%  '
$ws.Range("D25").Value = '2.061.76'
$ws.Range("E25").Value = '  -6.18%  '
$ws.Range("D26").NumberFormat = "@"
$ws.Range("D26").Value = '153.72'
$ws.Range("E26").Value = '  +0.14%  '
$ws.Range("D27").NumberFormat = "@"
$ws.Range("D27").Value = '19.09'
$ws.Range("E27").Value = '  -1.69%  '
$ws.Range("D28").NumberFormat = "@"
$ws.Range("D28").Value = '5.658'
$ws.Range("E28").Value = '  -7.53%  '
$ws.Range("D29").NumberFormat = "@"
$ws.Range("D29").Value = '1.955'
$ws.Range("E29").Value = '  -2.89%  '
$ws.Range("D30").NumberFormat = "@"
$ws.Range("D30").Value = '117.37'
$ws.Range("E30").Value = '  -0.13%  '
$ws.Range("D31").NumberFormat = "@"
$ws.Range("D31").Value = '0.9339'
$ws.Range("E31").Value = '  -4.88%  '
$ws.Range("D32").NumberFormat = "@"
$ws.Range("D32").Value = '0.09248'
$ws.Range("E32").Value = '  -1.79%  '
$ws.Range("D33").NumberFormat = "@"
$ws.Range("D33").Value = '5.270'
$ws.Range("E33").Value = '  -1.90%  '
$ws.Range("D34").NumberFormat = "@"
$ws.Range("D34").Value = '1.318'
$ws.Range("E34").Value = '  -2.21%  '
$ws.Range("D35").NumberFormat = "@"
$ws.Range("D35").Value = '3.292'
$ws.Range("E35").Value = '  -5.63%  '
$ws.Range("D36").NumberFormat = "@"
$ws.Range("D36").Value = '0.05845'
$ws.Range("E36").Value = '  -5.06%  '
$ws.Range("D37").NumberFormat = "@"
$ws.Range("D37").Value = '0.02121'
$ws.Range("E37").Value = '  -3.79%  '
$ws.Range("D38").NumberFormat = "@"
$ws.Range("D38").Value = '1.139'
$ws.Range("E38").Value = '  -1.93%  '
$ws.Range("D39").NumberFormat = "@"
$ws.Range("D39").Value = '7.810'
$ws.Range("E39").Value = '  +2.60%  '
$ws.Range("D40").NumberFormat = "@"
$ws.Range("D40").Value = '0.5580'
$ws.Range("E40").Value = '  -2.24%  '
$ws.Range("D41").NumberFormat = "@"
$ws.Range("D41").Value = '9.865'
$ws.Range("E41").Value = '  -2.83%  '
$ws.Range("D42").NumberFormat = "@"
$ws.Range("D42").Value = '0.1759'
$ws.Range("E42").Value = '  -1.99%  '
$ws.Range("D43").NumberFormat = "@"
$ws.Range("D43").Value = '11.66'
$ws.Range("E43").Value = '  -1.20%  '
$ws.Range("D44").NumberFormat = "@"
$ws.Range("D44").Value = '0.5253'
$ws.Range("E44").Value = '  -2.72%  '
$ws.Range("D45").NumberFormat = "@"
$ws.Range("D45").Value = '0.07008'
$ws.Range("E45").Value = '  -2.11%  '
$ws.Range("D46").NumberFormat = "@"
$ws.Range("D46").Value = '2.120'
$ws.Range("E46").Value = '  -10.51%  '
$ws.Range("B47").Value = 'NEARProtocol'
$ws.Range("C47").Value = 'https://coinranking.com/coin/DCrsaMv68+nearprotocol-near'
$ws.Range("D47").NumberFormat = "@"
$ws.Range("D47").Value = '1.826'
$ws.Range("E47").Value = '  -4.32%  '
$ws.Range("B48").Value = 'WEMIXToken'
$ws.Range("C48").Value = 'https://coinranking.com/coin/08CsQa-Ov+wemixtoken-wemix'
$ws.Range("D48").NumberFormat = "@"
$ws.Range("D48").Value = '1.109'
$ws.Range("E48").Value = '  -11.35%  '
$ws.Range("D49").NumberFormat = "@"
$ws.Range("D49").Value = '112.60'
$ws.Range("E49").Value = '  -1.37%  '
$ws.Range("E50").Value = '  -0.13%  '
$ws.Range("D51").NumberFormat = "@"
$ws.Range("D51").Value = '2.322'
$ws.Range("E51").Value = '  +0.22%  '
